$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New horizontally + vertically merged header block (rows 5-6) ---
# Row 5 = top header row (group headers + single-row headers that will span both rows)
$ws.Range("A5").Value = "simpleHeader1"
$ws.Range("B5").Value = "group1"
$ws.Range("D5").Value = "group2"
$ws.Range("F5").Value = "regexHeader"
$ws.Range("G5").Value = "regexHeader2"

# Row 6 = sub-header row, nested under the group headers above
$ws.Range("B6").Value = "mergedHeader1"
$ws.Range("D6").Value = "mergedHeader2"

# Apply the same centered style used by the existing merged header cells (C1:D1 / E1:F1)
# before merging, so every underlying cell keeps a single shared style.
$ws.Range("A5:G6").HorizontalAlignment = -4108

# Merge the header groups: vertical spans for single headers, horizontal spans for groups
$ws.Range("A5:A6").Merge()
$ws.Range("F5:F6").Merge()
$ws.Range("G5:G6").Merge()
$ws.Range("B6:C6").Merge()
$ws.Range("D6:E6").Merge()
$ws.Range("B5:C5").Merge()
$ws.Range("D5:E5").Merge()

# --- Extend data row 2 with matching blank formatted cells for the new columns ---
$ws.Range("D2").NumberFormat = "0"
$ws.Range("F2").NumberFormat = "0"
$ws.Range("G2").NumberFormat = "0"

# --- Marker cell used while authoring / selecting the new layout ---
$ws.Range("A10").WrapText = $False

# Final selection, matching where the user's cursor ended up
$ws.Range("E9").Select()
